$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$a2 = @"
Pay now
CADILLAC ESCALADE, 2023, Blue
N
85540
Date and Time of Issuing The Fine:
14 Jul 2025, 12:10 am
Location:
Ras Al khour St
Source:
Dubai Police
Amount:
AED 600
Payable Black Points:
-
Online declaration:
NO
Fine Number:
7037866556
Details:
Exceeding maximum speed limit by not more than 30 km h
Dispute:
Please contact Dubai Police for details about disputing your fine.
"@

$a3 = @"
Pay now
KIA K5, 2023, Black
DD
81392
Date and Time of Issuing The Fine:
11 Jul 2025, 8:30 am
Location:
Dubai Alain Road
Source:
Dubai Police
Amount:
AED 600
Payable Black Points:
-
Online declaration:
NO
Fine Number:
7037841032
Details:
Exceeding maximum speed limit by not more than 30 km h
Dispute:
Please contact Dubai Police for details about disputing your fine.
"@

# Trim the trailing newline that the here-string introduces before the closing tag.
$a2 = $a2.TrimEnd("`r", "`n")
$a3 = $a3.TrimEnd("`r", "`n")

$ws.Range("A2").Value = $a2
$ws.Range("A3").Value = $a3

$ws.Range("A5").EntireRow.Delete()
$ws.Range("A4").EntireRow.Delete()
